# Update gh-pages output data (view/like counts) for two sheets:
# "展览" (Exhibitions) and "全部类型" (All types)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5544
$ws1.Range("F9").Value = 525
$ws1.Range("F10").Value = 16

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5544
$ws4.Range("F11").Value = 525
$ws4.Range("F12").Value = 16
